$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B12").Value = 3
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 33.33
$ws.Range("E12").Value = 298.9701

$ws.Range("A13").Value = 13
$ws.Range("B13").Value = 3
$ws.Range("C13").Value = 2
$ws.Range("D13").Value = 33.33
$ws.Range("E13").Value = 298.9701

$ws.Range("A14").Value = 14
$ws.Range("B14").Value = 3
$ws.Range("C14").Value = 4
$ws.Range("D14").Value = 33.33
$ws.Range("E14").Value = 298.9701
